$wb = $excel.ActiveWorkbook

# ---- YDS sheet: append new game play-by-play yardage logs ----
$ydsWs = $wb.Worksheets.Item("YDS")
$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 4 2 5 4 3 0 7 5 2 1 10 6 5 3 1 3 3 11 9 1 27 7 3 5 1 3 1 4 3 2 2 23 4 -1 5 2 -2 4 11 3 2 0"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 9 4 5 13 7 9 7 17 4 3 5 9 2 0 15"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 3 6 7 -1 8 4 15 -1 3 5 1 7 12 9 2 4 6 19 0 4 4 6 5 -1 2 4 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 2 5 35 7 5 13 23 17 6 2 28 6 32 5 1 18 39 11 8"

# ---- ST sheet: append new game special-teams play logs ----
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 63"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 24"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 20 18"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 44 63 57 45"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 11 0 9 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 0 0"

# ---- OFF sheet: updated season totals ----
$wb.Worksheets.Item("OFF").Range("C2").Value = 378
$wb.Worksheets.Item("OFF").Range("D2").Value = 27
$wb.Worksheets.Item("OFF").Range("E2").Value = 20
$wb.Worksheets.Item("OFF").Range("F2").Value = 133
$wb.Worksheets.Item("OFF").Range("G2").Value = 117
$wb.Worksheets.Item("OFF").Range("I2").Value = 19
$wb.Worksheets.Item("OFF").Range("J2").Value = 70
$wb.Worksheets.Item("OFF").Range("L2").Value = 597
$wb.Worksheets.Item("OFF").Range("M2").Value = 393
$wb.Worksheets.Item("OFF").Range("O2").Value = 32
$wb.Worksheets.Item("OFF").Range("Q2").Value = 1112
$wb.Worksheets.Item("OFF").Range("B3").Value = 21
$wb.Worksheets.Item("OFF").Range("C3").Value = 411
$wb.Worksheets.Item("OFF").Range("F3").Value = 222
$wb.Worksheets.Item("OFF").Range("G3").Value = 93
$wb.Worksheets.Item("OFF").Range("H3").Value = 49
$wb.Worksheets.Item("OFF").Range("I3").Value = 111
$wb.Worksheets.Item("OFF").Range("J3").Value = 135
$wb.Worksheets.Item("OFF").Range("N3").Value = 37

# ---- DEF sheet: updated season totals ----
$wb.Worksheets.Item("DEF").Range("C2").Value = 400
$wb.Worksheets.Item("DEF").Range("D2").Value = 20
$wb.Worksheets.Item("DEF").Range("F2").Value = 120
$wb.Worksheets.Item("DEF").Range("G2").Value = 98
$wb.Worksheets.Item("DEF").Range("J2").Value = 49
$wb.Worksheets.Item("DEF").Range("L2").Value = 606
$wb.Worksheets.Item("DEF").Range("M2").Value = 362
$wb.Worksheets.Item("DEF").Range("Q2").Value = 1053
$wb.Worksheets.Item("DEF").Range("C3").Value = 389
$wb.Worksheets.Item("DEF").Range("D3").Value = 8
$wb.Worksheets.Item("DEF").Range("E3").Value = 74
$wb.Worksheets.Item("DEF").Range("F3").Value = 234
$wb.Worksheets.Item("DEF").Range("G3").Value = 67
$wb.Worksheets.Item("DEF").Range("H3").Value = 65
$wb.Worksheets.Item("DEF").Range("I3").Value = 122
$wb.Worksheets.Item("DEF").Range("J3").Value = 105
$wb.Worksheets.Item("DEF").Range("N3").Value = 52

# ---- ST sheet: updated season totals ----
$stWs.Range("B2").Value = 163
$stWs.Range("D2").Value = 144
$stWs.Range("F2").Value = 137
$stWs.Range("G2").Value = 134
$stWs.Range("J2").Value = 57
$stWs.Range("K2").Value = 54
$stWs.Range("L2").Value = 42
$stWs.Range("M2").Value = 32
$stWs.Range("B3").Value = 120

# ---- TURNS sheet: updated season totals ----
$wb.Worksheets.Item("TURNS").Range("C2").Value = 21
$wb.Worksheets.Item("TURNS").Range("E2").Value = 17
$wb.Worksheets.Item("TURNS").Range("D3").Value = 16
$wb.Worksheets.Item("TURNS").Range("E3").Value = 16

# ---- PEN sheet: updated season totals ----
$wb.Worksheets.Item("PEN").Range("B2").Value = 33
$wb.Worksheets.Item("PEN").Range("D4").Value = 13
